# Fix Training Data Issue (#48)
# The "Date" column (BF) held strings like "6-29-2013-14" which mixes the
# game date with the season label. The NBA stats site showed the date one
# day off, so the real game date is 2014-06-29. Normalize BF2:BF31 from
# "6-29-2013-14" to the correct ISO-style date string "2014-06-29".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 31
$dateCol = 58  # column BF

$oldText = "6-29-2013-14"
$newText = "2014-06-29"

# Force the target cells to text format first so the ISO-looking
# replacement string isn't silently reinterpreted as a date serial.
$ws.Range("BF$firstRow`:BF$lastRow").NumberFormat = "@"

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $dateCol)
    if ($cell.Value() -eq $oldText) {
        $cell.Value = $newText
    }
}
